# Apply the "add 2022-Q3 data" change:
#  1. Insert a new row (for 2022-Q3) at the top of the data in the "总计"
#     (totals) sheet, pushing the existing quarters down by one row.
#  2. Insert a brand-new worksheet named "2022-Q3" right after "总计",
#     populated with the per-fund holding data for that quarter.
#
# All the other existing quarter sheets (2022-Q2 .. 2020-Q4) are left
# untouched content-wise; they simply shift one position to the right as a
# natural consequence of inserting the new sheet after "总计".

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force the cell to store a text value even when it looks numeric
    # (e.g. "32.65"), matching the workbook's convention of keeping these
    # figures as text, and then strip the temporary "@" number format so
    # the cell is left with no explicit style, matching the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Part 1: "总计" sheet - insert the 2022-Q3 summary row at the top of the
# data (row 2), shifting the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$existingRows = @()
$r = 2
while ($true) {
    $dateCell = $totalSheet.Cells.Item($r, 2)
    $dateVal = $dateCell.Value()
    if ($dateVal -eq $null -or $dateVal -eq "") {
        break
    }
    $countVal = $totalSheet.Cells.Item($r, 3).Value()
    $valueVal = $totalSheet.Cells.Item($r, 4).Value()
    $existingRows += , @($dateVal, $countVal, $valueVal)
    $r = $r + 1
}

# Write rows back out shifted down by one, with the new 2022-Q3 row first.
$newRows = @(, @("2022-Q3", 22, 3.93))
$newRows += $existingRows

# Column A holds a plain numeric 0-based row index, styled (bold /
# centered / bordered) the same way as the header row. Use the existing
# A2 cell (style "s=2") as the formatting template for every A-column
# cell, copying it down before writing the value.
$aStyleSource = $totalSheet.Cells.Item(2, 1)

$rowIndex = 2
foreach ($row in $newRows) {
    $aCell = $totalSheet.Cells.Item($rowIndex, 1)
    $aStyleSource.Copy($aCell)
    $aCell.Value = ($rowIndex - 2)

    $totalSheet.Cells.Item($rowIndex, 2).Value = $row[0]
    $totalSheet.Cells.Item($rowIndex, 3).Value = $row[1]
    $totalSheet.Cells.Item($rowIndex, 4).Value = $row[2]
    $rowIndex = $rowIndex + 1
}

# ---------------------------------------------------------------------
# Part 2: insert the new "2022-Q3" worksheet right after "总计". Cross-
# sheet Range.Copy does not carry formatting in this host, so instead
# duplicate the existing "2022-Q2" sheet (a full Worksheet.Copy keeps all
# styles/column widths/etc.) and overwrite its data in place - the
# header row (row 1) is already identical text, so it is left untouched.
# ---------------------------------------------------------------------
$afterQ2 = $wb.Worksheets.Item("2022-Q2")
$afterQ2.Copy($afterQ2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Style template for the column-A numeric index, taken from the new
# sheet's own (already duplicated) A2 cell.
$indexStyleSource = $q3.Cells.Item(2, 1)

$fundRows = @(
    @("050026", "博时医疗保健行业混合A", "32.65", "88.15", "2.98", "0.9730", 9),
    @("009630", "浦银安盛ESG责任投资混合A", "8.65", "84.74", "8.26", "0.7145", 2),
    @("001230", "鹏华医药科技股票", "15.31", "80.22", "3.73", "0.5711", 6),
    @("009631", "浦银安盛ESG责任投资混合C", "4.98", "84.74", "8.26", "0.4113", 2),
    @("001766", "上投摩根医疗健康股票A", "8.99", "82.99", "2.91", "0.2616", 10),
    @("007067", "浦银安盛先进制造混合C", "2.06", "77.53", "8.15", "0.1679", 1),
    @("000968", "广发中证养老产业指数A", "10.84", "93.98", "1.38", "0.1496", 9),
    @("012202", "中加消费优选混合A", "3.52", "73.15", "3.97", "0.1397", 5),
    @("003581", "新疆前海联合国民健康产业灵活配置混合A", "2.10", "93.17", "6.45", "0.1354", 5),
    @("001551", "天弘中证医药100指数型发起式 C", "8.58", "95.24", "1.37", "0.1175", 8),
    @("007066", "浦银安盛先进制造混合A", "1.27", "77.53", "8.15", "0.1035", 1),
    @("001550", "天弘中证医药100指数型发起式 A", "5.31", "95.24", "1.37", "0.0727", 8),
    @("004050", "华夏新锦升灵活配置混合A", "1.22", "64.62", "3.27", "0.0399", 9),
    @("012203", "中加消费优选混合C", "0.57", "73.15", "3.97", "0.0226", 5),
    @("007111", "新疆前海联合国民健康产业灵活配置混合C", "0.23", "93.17", "6.45", "0.0148", 5),
    @("002982", "广发中证养老产业指数C", "0.92", "93.98", "1.38", "0.0127", 9),
    @("516560", "华宝养老ETF", "0.74", "98.01", "1.44", "0.0107", 9),
    @("011895", "博时医疗保健行业混合C", "0.13", "88.15", "2.98", "0.0039", 9),
    @("004051", "华夏新锦升灵活配置混合C", "0.06", "64.62", "3.27", "0.0020", 9),
    @("014932", "上投摩根医疗健康股票C", "0.04", "82.99", "2.91", "0.0012", 10),
    @("014328", "格林新兴产业混合C", "0.03", "64.92", "3.24", "0.0010", 10),
    @("014327", "格林新兴产业混合A", "0.01", "64.92", "3.24", "0.0003", 10)
)

$rowIndex = 2
foreach ($fund in $fundRows) {
    $aCell = $q3.Cells.Item($rowIndex, 1)
    $indexStyleSource.Copy($aCell)
    $aCell.Value = $rowIndex - 2

    Set-TextValue $q3.Cells.Item($rowIndex, 2) $fund[0]
    Set-TextValue $q3.Cells.Item($rowIndex, 3) $fund[1]
    Set-TextValue $q3.Cells.Item($rowIndex, 4) $fund[2]
    Set-TextValue $q3.Cells.Item($rowIndex, 5) $fund[3]
    Set-TextValue $q3.Cells.Item($rowIndex, 6) $fund[4]
    Set-TextValue $q3.Cells.Item($rowIndex, 7) $fund[5]
    $q3.Cells.Item($rowIndex, 8).Value = $fund[6]

    $rowIndex = $rowIndex + 1
}

$q3.Range("A1").Select()
